$d = $word.ActiveDocument

# The document has one section whose primary header/footer and
# even-page header/footer each carry one inline picture:
#   - Headers: the BTec "BTec_Logo-Orange" logo (currently named image1.jpg)
#   - Footers: the Pearson "PearsonLogo.png" logo (currently named image2.png)
#
# This commit swaps the picture names: the BTec logos become "image2.jpg"
# and the Pearson logos become "image1.png".

$sec = $d.Sections.Item(1)

$hdr1 = $sec.Headers.Item(1)
$hdr2 = $sec.Headers.Item(2)
$ftr1 = $sec.Footers.Item(1)
$ftr2 = $sec.Footers.Item(2)

# Headers: renaming the inline picture directly works fine.
$hdr1.Range.InlineShapes.Item(1).Name = "image2.jpg"
$hdr2.Range.InlineShapes.Item(1).Name = "image2.jpg"

# Footers: renaming InlineShapes.Item(1) directly on a footer-owned
# picture trips this host's "stale handle" guard, so round-trip the
# picture through ConvertToShape/ConvertToInlineShape (which re-anchors
# it) to apply the rename, then flip it back to an inline picture.
$ftr1Pic = $ftr1.Range.InlineShapes.Item(1).ConvertToShape()
$ftr1Pic.Name = "image1.png"
$ftr1Pic.ConvertToInlineShape() | Out-Null

$ftr2Pic = $ftr2.Range.InlineShapes.Item(1).ConvertToShape()
$ftr2Pic.Name = "image1.png"
$ftr2Pic.ConvertToInlineShape() | Out-Null

Write-Output "Renamed header/footer logo pictures"
